$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("formulaire")
Write-Host $ws.Name
